# Rows 6 and 7 of the "Artfynd" sheet hold two observation records that were
# re-sorted: the record previously on row 6 now belongs on row 7, and vice
# versa. Swap the per-record fields between the two rows; columns that are
# identical for both records (P, S, T, U, V, W, Y, AA, AD, AE, AG, AT, AW,
# AX, AY, ...) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $addr6 = "$col" + "6"
    $addr7 = "$col" + "7"
    $v6 = $ws.Range($addr6).Value2
    $v7 = $ws.Range($addr7).Value2
    $ws.Range($addr6).Value2 = $v7
    $ws.Range($addr7).Value2 = $v6
}

# The (blank) "Bestämningsmetod" marker cell moves from AF6 to AF7 along
# with the rest of the row-7-bound record. Recreate it as an explicit empty
# text cell (matching the other blank-but-present text cells on the row,
# e.g. I7/AT7/AY7) rather than merely leaving it unset/absent.
$ws.Range("AF6").ClearContents()
$ws.Range("AF7").Value = "'"
$ws.Range("AF7").ClearFormats()
